$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Douwe")

# Fill in newly logged work session (row 5 description + row 6 new entry)
$ws2.Range("D5").Value = "Verder gewerkt aan de code van de eerste kamer"

$ws2.Range("A6").Value = 0.4375
$ws2.Range("B6").Value = 0.5625
$ws2.Range("D6").Value = "Begin gemaakt aan de derde ruimte van onze game. Ook gewerkt aan de geheime ruimte"

# Make the "Douwe" sheet the active/selected one, with D6 as the active cell
$ws2.Activate()
$null = $ws2.Range("D6").Select()
